# Apply cryptos list price/volume update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.434.43"
$ws.Range("E2").Value = "  +3.43%  "
$ws.Range("D3").Value = "2.310.96"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'311.33"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").Value = "'103.84"
$ws.Range("E6").Value = "  +7.60%  "
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("E9").Value = "  +8.71%  "
$ws.Range("D10").Value = "'36.45"
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("D11").Value = "'0.0814"
$ws.Range("E11").Value = "  +3.78%  "
$ws.Range("D12").Value = "'51.70"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "'7.06"
$ws.Range("E14").Value = "  +3.85%  "
$ws.Range("D15").Value = "2.669.38"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("D16").Value = "'15.07"
$ws.Range("E16").Value = "  +3.49%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'0.811"
$ws.Range("E17").Value = "  +3.26%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.220.28"
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("D19").Value = "43.329.68"
$ws.Range("E19").Value = "  +3.56%  "
$ws.Range("D20").Value = "'12.35"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "0.0₃0931"
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("D22").Value = "'6.19"
$ws.Range("E22").Value = "  +3.99%  "
$ws.Range("D23").Value = "'68.22"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").Value = "'242.58"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("E26").Value = "  +2.55%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "'24.88"
$ws.Range("E28").Value = "  +6.20%  "
$ws.Range("E29").Value = "  +8.07%  "
$ws.Range("D30").Value = "'36.88"
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("D31").Value = "'9.69"
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("D32").Value = "'168.22"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").Value = "'5.29"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'18.08"
$ws.Range("E35").Value = "  +3.67%  "
$ws.Range("E36").Value = "  +6.17%  "
$ws.Range("D37").Value = "'0.0743"
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("D38").Value = "'3.07"
$ws.Range("E38").Value = "  -2.52%  "
$ws.Range("E39").Value = "  +4.77%  "
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("E41").Value = "  +8.13%  "
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("D43").Value = "'2.50"
$ws.Range("E43").Value = "  +8.76%  "
$ws.Range("E44").Value = "  +6.24%  "
$ws.Range("D45").Value = "1.985.57"
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("D46").Value = "'19.10"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("E47").Value = "  +4.87%  "
$ws.Range("D48").Value = "'9.97"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").Value = "'55.94"
$ws.Range("E49").Value = "  +5.56%  "
$ws.Range("D50").Value = "'2.95"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "'1.59"
$ws.Range("E51").Value = "  +9.68%  "
